$d = $word.ActiveDocument
$d.Content.Find.Execute("^lAnggota / NIDN <<", $false, $false, $false, $false, $false, $true, 1, $false, "Anggota / NIDN <<", 2)
